# Auto-generated Excel COM-interop script to apply cryptos.xlsx update
# Commit: "Updated cryptos list on Thu May 16 11:36:15 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: address -> new text value. NumberFormat "@" is applied first
# so that values which look numeric (e.g. "163.22", "0.0000104") are stored as
# exact text, matching the original inline-string cells instead of being coerced
# into floating point numbers (which would silently lose formatting like trailing
# zeros, e.g. "163.00" -> 163).
$updates = @(
    @("D2", "66.187.75"),
    @("E2", "  +6.01%  "),
    @("D3", "2.999.32"),
    @("E3", "  +3.57%  "),
    @("E4", "  +0.01%  "),
    @("D5", "581.68"),
    @("E5", "  +2.66%  "),
    @("D6", "163.22"),
    @("E6", "  +13.64%  "),
    @("E7", "  -0.03%  "),
    @("D9", "2.996.70"),
    @("E9", "  +3.52%  "),
    @("D10", "6.53"),
    @("E10", "  -4.96%  "),
    @("D11", "0.154"),
    @("E11", "  +3.90%  "),
    @("D12", "0.455"),
    @("E12", "  +5.95%  "),
    @("E13", "  +6.15%  "),
    @("D14", "34.60"),
    @("E14", "  +5.70%  "),
    @("E15", "  -0.62%  "),
    @("D16", "66.183.46"),
    @("E16", "  +6.04%  "),
    @("D17", "3.497.04"),
    @("E17", "  +3.45%  "),
    @("D18", "6.90"),
    @("E18", "  +5.08%  "),
    @("D19", "2.997.95"),
    @("E19", "  +3.43%  "),
    @("D20", "453.12"),
    @("E20", "  +6.67%  "),
    @("E21", "  +6.30%  "),
    @("D22", "0.685"),
    @("E22", "  +4.33%  "),
    @("D23", "7.34"),
    @("E23", "  +7.20%  "),
    @("D24", "82.31"),
    @("E24", "  +4.78%  "),
    @("E25", "  +14.59%  "),
    @("D26", "12.29"),
    @("E26", "  +3.69%  "),
    @("D27", "10.39"),
    @("E27", "  +3.61%  "),
    @("E28", "  +0.04%  "),
    @("D29", "8.13"),
    @("E29", "  +14.73%  "),
    @("E30", "  +19.48%  "),
    @("B31", "PancakeSwap"),
    @("C31", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"),
    @("D31", "2.62"),
    @("E31", "  +5.80%  "),
    @("B32", "PEPE"),
    @("C32", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"),
    @("D32", "0.0000104"),
    @("E32", "  -4.82%  "),
    @("D33", "27.23"),
    @("E33", "  +6.13%  "),
    @("D34", "0.110"),
    @("E34", "  +5.22%  "),
    @("E35", "  +0.00%  "),
    @("D36", "0.991"),
    @("E36", "  +4.98%  "),
    @("E37", "  +8.21%  "),
    @("D38", "2.06"),
    @("E38", "  +9.23%  "),
    @("D39", "49.48"),
    @("E39", "  +2.07%  "),
    @("E40", "  +1.58%  "),
    @("D41", "0.311"),
    @("E41", "  +17.12%  "),
    @("D42", "44.21"),
    @("E42", "  +7.42%  "),
    @("E43", "  +7.20%  "),
    @("D44", "8.45"),
    @("E44", "  +5.36%  "),
    @("D45", "401.28"),
    @("E45", "  +13.09%  "),
    @("D46", "0.0358"),
    @("E46", "  +6.33%  "),
    @("D47", "2.769.45"),
    @("E47", "  +2.03%  "),
    @("D48", "133.47"),
    @("E48", "  +0.22%  "),
    @("D50", "23.87"),
    @("E50", "  +12.70%  "),
    @("E51", "  +4.11%  ")
)

foreach ($u in $updates) {
    $addr = $u[0]
    $val = $u[1]
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}
